$d = $word.ActiveDocument
$d.Content.Find.Execute("Scripts SQL de création", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Scripts Postgresql de création", 2)
